# Regenerate all penyata to follow new data and format
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Text relabelling: "Penandaan Fail" sub-items gain a "Semakan" prefix.
#    The same four labels are reused in three different sections further
#    down the sheet (Penandaan Fail / Laporan Atas Talian / JPPM-JDM-JDRM),
#    so every occurrence needs updating.
# ---------------------------------------------------------------------------
$ws.Range("C16").Value = "Semakan Kali Pertama"
$ws.Range("C17").Value = "Semakan Kali Kedua"
$ws.Range("C18").Value = "Semakan Kali Ketiga"
$ws.Range("C19").Value = "Semakan Kali Keempat"

$ws.Range("C22").Value = "Semakan Kali Pertama"
$ws.Range("C23").Value = "Semakan Kali Kedua"
$ws.Range("C24").Value = "Semakan Kali Ketiga"
$ws.Range("C25").Value = "Semakan Kali Keempat"

$ws.Range("C28").Value = "Semakan Kali Pertama"
$ws.Range("C29").Value = "Semakan Kali Kedua"
$ws.Range("C30").Value = "Semakan Kali Ketiga"
$ws.Range("C31").Value = "Semakan Kali Keempat"

# ---------------------------------------------------------------------------
# 2. Competition names switch from ALL CAPS to Title Case.
# ---------------------------------------------------------------------------
$ws.Range("C34").Value = "Bouquet Kreatif"
$ws.Range("C35").Value = "Kad Raya Untuk Guruku"
$ws.Range("C36").Value = "Riang Ria Kuih Raya"
$ws.Range("C37").Value = "Creative Collage"

# ---------------------------------------------------------------------------
# 3. Updated "Penandaan Fail" actuals for "Kali Ketiga" (row 18).
# ---------------------------------------------------------------------------
$ws.Range("D18").Value = 3606
$ws.Range("E18").Value = 2370

# ---------------------------------------------------------------------------
# 4. Three new competitions added under "Penyertaan Pertandingan".
# ---------------------------------------------------------------------------
$ws.Range("C38").Value = "Lompat Getah"
$ws.Range("D38").Value = 500
$ws.Range("E38").Value = 0

$ws.Range("C39").Value = "Theme Party"
$ws.Range("D39").Value = 100
$ws.Range("E39").Value = 0

$ws.Range("C40").Value = "Hari Koperasi"
$ws.Range("D40").Value = 0
$ws.Range("E40").Value = 0

# ---------------------------------------------------------------------------
# 5. Merged-cell layout tweaks to match the reflowed header / section rows.
# ---------------------------------------------------------------------------
$ws.Range("D4:G4").Merge()
$ws.Range("B5:C5").Merge()
$ws.Range("B12:F12").Merge()

$ws.Range("B15:C15").UnMerge()

$ws.Range("B21:C21").UnMerge()
$ws.Range("B21:E21").Merge()

$ws.Range("B27:C27").UnMerge()
$ws.Range("B27:E27").Merge()

$ws.Range("B33:C33").UnMerge()
$ws.Range("B33:E33").Merge()

$ws.Range("B43:E43").Merge()

# ---------------------------------------------------------------------------
# 6. Logo picture repositioned / resized.
# ---------------------------------------------------------------------------
$shp = $ws.Shapes.Item(1)
$shp.Left = 41.2125
$shp.Top = 14.25
$shp.Width = 46.5
$shp.Height = 47.25

# ---------------------------------------------------------------------------
# 7. Drop the trailing blank formatted row (1001) so the sheet ends at 1000.
# ---------------------------------------------------------------------------
$ws.Rows.Item(1001).EntireRow.Delete()

# ---------------------------------------------------------------------------
# 8. Page setup: fit to one page tall/wide, centre horizontally, no header
#    / footer margin reservation.
# ---------------------------------------------------------------------------
$ws.PageSetup.FitToPagesWide = 1
$ws.PageSetup.FitToPagesTall = 1
$ws.PageSetup.CenterHorizontally = $true
$ws.PageSetup.HeaderMargin = 0
$ws.PageSetup.FooterMargin = 0
